# Remove the leftover empty "AutoShape" placeholders (Google image-search
# result artifacts) that sit behind the real pictures on several slides.
# Each entry maps a 1-based slide index to the set of shape Ids (p:cNvPr id)
# that must be deleted from that slide.

$p = $ppt.ActivePresentation

$toDelete = @{
    1 = @(20482)
    2 = @(20482, 21506)
    3 = @(1029, 1031)
    5 = @(5124, 5126, 5134, 5136, 5138, 5141, 5143)
    8 = @(19458, 19460, 19462)
}

foreach ($slideIndex in $toDelete.Keys) {
    $s = $p.Slides.Item($slideIndex)
    $ids = $toDelete[$slideIndex]
    for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
        $sh = $s.Shapes.Item($i)
        if ($ids -contains $sh.Id) {
            $sh.Delete()
        }
    }
}
